$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1463.5
$ws.Range("I40").Value = 1437.125
$ws.Range("J40").Value = 1498.6666
$ws.Range("K40").Value = 1437.125
$ws.Range("L40").Value = 1498.6666
$ws.Range("M40").Value = -1262.125
$ws.Range("N40").Value = -1848.6666
$ws.Range("H88").Value = 26866.75
$ws.Range("I88").Value = 1350.75
$ws.Range("J88").Value = 39624.75
$ws.Range("K88").Value = 1350.75
$ws.Range("L88").Value = 39624.75
$ws.Range("M88").Value = -944.75
$ws.Range("N88").Value = -40436.75
$ws.Range("H91").Value = 26866.75
$ws.Range("I91").Value = 1350.75
$ws.Range("J91").Value = 39624.75
$ws.Range("K91").Value = 1350.75
$ws.Range("L91").Value = 39624.75
$ws.Range("M91").Value = 53.25
$ws.Range("N91").Value = -42432.75
$ws.Range("H132").Value = 9412.111000000001
$ws.Range("I132").Value = 11835
$ws.Range("J132").Value = 4566.3335
$ws.Range("K132").Value = 35505
$ws.Range("L132").Value = 13699.0005
$ws.Range("M132").Value = -32975
$ws.Range("N132").Value = -18759.0005
$ws.Range("H133").Value = 52197.5
$ws.Range("J133").Value = 52197.5
$ws.Range("L133").Value = 52197.5
$ws.Range("N133").Value = -62317.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2019314.2
$ws.Range("I32").Value = 2533766.8
$ws.Range("K32").Value = 2533766.8
$ws.Range("M32").Value = -2533479.8
$ws.Range("H122").Value = 51234.5
$ws.Range("I122").Value = 51234.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 153703.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -151253.5
$ws.Range("N122").Value = $null
$ws.Range("H123").Value = 98429
$ws.Range("J123").Value = 98429
$ws.Range("L123").Value = 98429
$ws.Range("N123").Value = -108229
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H135").Value = 19321.75
$ws.Range("J135").Value = 19321.75
$ws.Range("L135").Value = 19321.75
$ws.Range("N135").Value = -29461.75
$ws.Range("H139").Value = 61893.668
$ws.Range("J139").Value = 61893.668
$ws.Range("L139").Value = 61893.668
$ws.Range("N139").Value = -72173.66800000001
$ws.Range("H140").Value = 65724.39999999999
$ws.Range("J140").Value = 65724.39999999999
$ws.Range("L140").Value = 65724.39999999999
$ws.Range("N140").Value = -76084.39999999999

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 695.875
$ws.Range("I22").Value = 684.2857
$ws.Range("J22").Value = 777
$ws.Range("K22").Value = 684.2857
$ws.Range("L22").Value = 777
$ws.Range("M22").Value = -511.2857
$ws.Range("N22").Value = -1123
$ws.Range("H86").Value = 1991.9302
$ws.Range("I86").Value = 1780.742
$ws.Range("J86").Value = 2537.5
$ws.Range("K86").Value = 1780.742
$ws.Range("L86").Value = 2537.5
$ws.Range("M86").Value = -657.742
$ws.Range("N86").Value = -4783.5
$ws.Range("H89").Value = 1991.9302
$ws.Range("I89").Value = 1780.742
$ws.Range("J89").Value = 2537.5
$ws.Range("K89").Value = 8903.709999999999
$ws.Range("L89").Value = 12687.5
$ws.Range("M89").Value = -3287.709999999999
$ws.Range("N89").Value = -23919.5
$ws.Range("H132").Value = 64440
$ws.Range("J132").Value = 69253.336
$ws.Range("L132").Value = 69253.336
$ws.Range("N132").Value = -79373.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 1894523.4
$ws.Range("I107").Value = 3472518.2
$ws.Range("J107").Value = 929.5333000000001
$ws.Range("K107").Value = 3472518.2
$ws.Range("L107").Value = 929.5333000000001
$ws.Range("M107").Value = -3470598.2
$ws.Range("N107").Value = -4769.5333
$ws.Range("H122").Value = 1604.4073
$ws.Range("I122").Value = 1112.4166
$ws.Range("J122").Value = 1998
$ws.Range("K122").Value = 3337.2498
$ws.Range("L122").Value = 5994
$ws.Range("M122").Value = -887.2498000000001
$ws.Range("N122").Value = -10894
$ws.Range("H132").Value = 6175266.5
$ws.Range("I132").Value = 2285.2144
$ws.Range("J132").Value = 12823092
$ws.Range("K132").Value = 6855.6432
$ws.Range("L132").Value = 38469276
$ws.Range("M132").Value = -4325.6432
$ws.Range("N132").Value = -38474336
$ws.Range("H133").Value = 51000
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 51000
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 51000
$ws.Range("M133").Value = $null
$ws.Range("N133").Value = -56060
$ws.Range("H135").Value = 45845.43
$ws.Range("J135").Value = 45845.43
$ws.Range("L135").Value = 45845.43
$ws.Range("N135").Value = -55985.43

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 643.7143
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 643.7143
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 1931.1429
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -4427.1429
$ws.Range("H131").Value = 889.95123
$ws.Range("J131").Value = 1147.2693
$ws.Range("L131").Value = 3441.8079
$ws.Range("N131").Value = -13521.8079

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 8209.8125
$ws.Range("J122").Value = 2201.6
$ws.Range("L122").Value = 6604.799999999999
$ws.Range("N122").Value = -11504.8
$ws.Range("H134").Value = 33708.668
$ws.Range("J134").Value = 33708.668
$ws.Range("L134").Value = 101126.004
$ws.Range("N134").Value = -106196.004
$ws.Range("H141").Value = 65013.125
$ws.Range("J141").Value = 65013.125
$ws.Range("L141").Value = 65013.125
$ws.Range("N141").Value = -75373.125

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 12493.6
$ws.Range("I22").Value = 540.2
$ws.Range("J22").Value = 18470.3
$ws.Range("K22").Value = 540.2
$ws.Range("L22").Value = 18470.3
$ws.Range("M22").Value = -245.2
$ws.Range("N22").Value = -19060.3
$ws.Range("H27").Value = 12493.6
$ws.Range("I27").Value = 540.2
$ws.Range("J27").Value = 18470.3
$ws.Range("K27").Value = 540.2
$ws.Range("L27").Value = 18470.3
$ws.Range("M27").Value = -433.2
$ws.Range("N27").Value = -18684.3
$ws.Range("H122").Value = 4766.6113
$ws.Range("I122").Value = 3533.1667
$ws.Range("J122").Value = 5383.3335
$ws.Range("K122").Value = 10599.5001
$ws.Range("L122").Value = 16150.0005
$ws.Range("M122").Value = -8149.500100000001
$ws.Range("N122").Value = -21050.0005
$ws.Range("H132").Value = 3184.973
$ws.Range("I132").Value = 2760.818
$ws.Range("J132").Value = 3807.0667
$ws.Range("K132").Value = 8282.454000000002
$ws.Range("L132").Value = 11421.2001
$ws.Range("M132").Value = -5752.454000000002
$ws.Range("N132").Value = -16481.2001
$ws.Range("H135").Value = 20875
$ws.Range("J135").Value = 20875
$ws.Range("L135").Value = 20875
$ws.Range("N135").Value = -31015
$ws.Range("H136").Value = 4631213.5
$ws.Range("I136").Value = 1165.6154
$ws.Range("J136").Value = 16669337
$ws.Range("K136").Value = 3496.8462
$ws.Range("L136").Value = 50008011
$ws.Range("M136").Value = -946.8462
$ws.Range("N136").Value = -50013111
$ws.Range("H138").Value = 44800
$ws.Range("J138").Value = 44800
$ws.Range("L138").Value = 44800
$ws.Range("N138").Value = -55080
$ws.Range("H140").Value = 53614.6
$ws.Range("J140").Value = 53614.6
$ws.Range("L140").Value = 53614.6
$ws.Range("N140").Value = -63974.6
$ws.Range("H141").Value = 57099.145
$ws.Range("J141").Value = 57099.145
$ws.Range("L141").Value = 57099.145
$ws.Range("N141").Value = -67459.14499999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 104809.336
$ws.Range("J46").Value = 104809.336
$ws.Range("L46").Value = 104809.336
$ws.Range("N46").Value = -105271.336
$ws.Range("H122").Value = 1320.2632
$ws.Range("I122").Value = 1345.625
$ws.Range("J122").Value = 1185
$ws.Range("K122").Value = 4036.875
$ws.Range("L122").Value = 3555
$ws.Range("M122").Value = -1586.875
$ws.Range("N122").Value = -8455
$ws.Range("H124").Value = 47119.332
$ws.Range("J124").Value = 47119.332
$ws.Range("L124").Value = 47119.332
$ws.Range("N124").Value = -56939.332
$ws.Range("H133").Value = 50715
$ws.Range("J133").Value = 50715
$ws.Range("L133").Value = 50715
$ws.Range("N133").Value = -60835
$ws.Range("H134").Value = 104809.336
$ws.Range("J134").Value = 104809.336
$ws.Range("L134").Value = 314428.008
$ws.Range("N134").Value = -319498.008
$ws.Range("H135").Value = 63200
$ws.Range("J135").Value = 63200
$ws.Range("L135").Value = 63200
$ws.Range("N135").Value = -73340
$ws.Range("H140").Value = 28766
$ws.Range("I140").Value = 10000
$ws.Range("J140").Value = 33457.5
$ws.Range("K140").Value = 10000
$ws.Range("L140").Value = 33457.5
$ws.Range("M140").Value = -4820
$ws.Range("N140").Value = -43817.5
$ws.Range("H141").Value = 131790.56
$ws.Range("I141").Value = 12000
$ws.Range("J141").Value = 146764.38
$ws.Range("K141").Value = 12000
$ws.Range("L141").Value = 146764.38
$ws.Range("M141").Value = -6820
$ws.Range("N141").Value = -157124.38
